$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Clear the whole sheet so we can rebuild the small lookup table cleanly.
# ---------------------------------------------------------------------------
$ws.Cells.Clear()

# ---------------------------------------------------------------------------
# 2. Re-write the table, now starting in column A instead of column B, with
#    the header row moved from row 2 to row 1 and the trailing "Others" / "?"
#    row removed.
# ---------------------------------------------------------------------------

# Header row
$ws.Range("A1").Value = "Classification #3"
$ws.Range("B1").Value = "Why defect was not identified during testing"

# Data rows (label, description)
$data = @(
    @("Ambiguous Requirement [SA]", "Requirement given was tested but interpreted wrongly."),
    @("Ambiguous Requirement [BR]", "Business requirement given was open to interpretation."),
    @("Ambiguous Requirement [Software]", "Test case failed to account for software updates."),
    @("Missing Requirement [SA]", "This requirement was not given by the design team."),
    @("Missing Requirement [BR]", "This requirement was not specified in the business document."),
    @("Missing Requirement [Software]", "This software requirement was not included."),
    @("Missing Test Case", "Test case not performed."),
    @("Not Within Scope", "This is outside the scope of the project."),
    @("Test & Production Environment Difference", "Test case was not chosen based on the impact analysis.  "),
    @("Test & Production Operation Difference", "This issue is specific to production operations."),
    @("Time Constraint", "Issue was not tested due to time constraint.")
)

$r = 2
foreach ($pair in $data) {
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
    $r++
}

Write-Host "values written"
